$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.541.32'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.736.99'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.91'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4940'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2669'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06293'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.729.21'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07045'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.71'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.595'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6125'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.0000'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.537.56'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007344'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +6.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9998'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.955.74'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.705'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.251'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '140.08'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.46'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.421'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '108.07'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.764'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.045'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08076'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.715'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04595'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.34%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.009'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6365'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.8967'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.016'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.400'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.005'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.94'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.403'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3907'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.872'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1189'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05397'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.23%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.820'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.57'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.269'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.81'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.98%  '
